$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update changed numeric values in existing rows (2-18) ---
$ws.Range("B2").Value = 0.012564715786985
$ws.Range("C2").Value = 0.24401343485806
$ws.Range("D2").Value = -0.059294523365166
$ws.Range("E2").Value = -0.08883472622775901
$ws.Range("F2").Value = 0.027489156085751
$ws.Range("B3").Value = 0.119301216398655
$ws.Range("C3").Value = 1.09091870030191
$ws.Range("D3").Value = -0.18292015852187
$ws.Range("E3").Value = -0.815881595142061
$ws.Range("F3").Value = 0.054958688836974
$ws.Range("B4").Value = -0.520020543637096
$ws.Range("C4").Value = 4.54730092679259
$ws.Range("D4").Value = -0.921670690352811
$ws.Range("E4").Value = -2.29731569051892
$ws.Range("F4").Value = 0.211875224941396
$ws.Range("B5").Value = 0.058099719588104
$ws.Range("C5").Value = 0.421422962693021
$ws.Range("D5").Value = -0.119117107958588
$ws.Range("E5").Value = -0.197471584550116
$ws.Range("F5").Value = 0.041313950928128
$ws.Range("B6").Value = 0.038523484528755
$ws.Range("C6").Value = 0.35761340334831
$ws.Range("D6").Value = -0.093191649949571
$ws.Range("E6").Value = -0.161883510159085
$ws.Range("F6").Value = 0.035791456191681
$ws.Range("B7").Value = 0.04728946340456101
$ws.Range("C7").Value = 0.359367288783417
$ws.Range("D7").Value = -0.09398158286654
$ws.Range("E7").Value = -0.160260990063469
$ws.Range("F7").Value = 0.038608567131846
$ws.Range("B8").Value = 0.03257350023650801
$ws.Range("C8").Value = 0.377692285484843
$ws.Range("D8").Value = -0.08251287122122
$ws.Range("E8").Value = -0.148379674097089
$ws.Range("F8").Value = 0.045422318853138
$ws.Range("B9").Value = 0.165203584871899
$ws.Range("C9").Value = 0.451641477581985
$ws.Range("D9").Value = -0.07992368353207101
$ws.Range("E9").Value = -0.28401170258503
$ws.Range("F9").Value = 0.063758950884435
$ws.Range("B10").Value = 0.165203584871899
$ws.Range("C10").Value = 0.451641477581985
$ws.Range("D10").Value = -0.07992368353207101
$ws.Range("E10").Value = -0.28401170258503
$ws.Range("B11").Value = 0.010615554946815
$ws.Range("C11").Value = 0.368357206529457
$ws.Range("D11").Value = -0.09826196746764401
$ws.Range("E11").Value = -0.174548875700433
$ws.Range("F11").Value = 0.027156011019518
$ws.Range("B12").Value = -4.88625842284954
$ws.Range("C12").Value = 4.34851015618375
$ws.Range("D12").Value = 7.19205840154032
$ws.Range("E12").Value = 3.83777003140382
$ws.Range("F12").Value = 2.65051434446095
$ws.Range("B13").Value = -5.12940141887619
$ws.Range("C13").Value = 4.21901419075841
$ws.Range("D13").Value = 8.87974220640168
$ws.Range("E13").Value = 6.19657528429093
$ws.Range("F13").Value = 3.56980331327197
$ws.Range("B14").Value = 0.518329924020808
$ws.Range("C14").Value = 1.06888955872832
$ws.Range("D14").Value = 0.425351537721022
$ws.Range("E14").Value = 0.8954643317453991
$ws.Range("F14").Value = 0.727404358922172
$ws.Range("B15").Value = 0.06294001163187
$ws.Range("C15").Value = 0.5295329386748431
$ws.Range("D15").Value = -0.112946262041283
$ws.Range("E15").Value = -0.300465802198912
$ws.Range("F15").Value = 0.045635177299235
$ws.Range("B16").Value = -0.054795281405736
$ws.Range("C16").Value = 0.5835011904303731
$ws.Range("D16").Value = 0.041003940196902
$ws.Range("E16").Value = -0.201571913031992
$ws.Range("F16").Value = 0.09344036926013301
$ws.Range("B17").Value = 2.40135549965032
$ws.Range("C17").Value = 12.4525138471674
$ws.Range("D17").Value = 15.542948420251
$ws.Range("E17").Value = 11.8730955924038
$ws.Range("F17").Value = 10.5974196996706
$ws.Range("B18").Value = -4.93695099676872
$ws.Range("C18").Value = 4.34252976444416
$ws.Range("D18").Value = 7.20118431619643
$ws.Range("E18").Value = 3.80336213466803
$ws.Range("F18").Value = 2.63018764559824

# --- Style fixes: B4 loses yellow highlight; B13 gains yellow highlight ---
$ws.Range("B4").Style = "Normal"
$ws.Range("B13").Interior.Color = 65535

# --- Add new rows 19 and 20 (copy label style/format from row 18) ---
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A19").Value = "arima1_0_2"
$ws.Range("B19").Value = -4.19505348367604
$ws.Range("C19").Value = 4.91159475378929
$ws.Range("D19").Value = 7.50972418272067
$ws.Range("E19").Value = 4.07263115094629
$ws.Range("F19").Value = 3.10182882796812
$ws.Range("A18").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A20").Value = "sarima212_001"
$ws.Range("B20").Value = 2.3412683459149
$ws.Range("C20").Value = 12.3331045517984
$ws.Range("D20").Value = 15.4148542117711
$ws.Range("E20").Value = 11.7800852320783
$ws.Range("F20").Value = 3.10182882796812

$excel.CutCopyMode = 0
